# Generate Report for Handoff
# Replaces the old "png" sample fixtures (3dd3ed8b..., 460c6bab..., cce5b3b7...)
# with the new md-based fixtures (calleeMd1.md, calleeMd2.md, callerMd1.md,
# callerMd2.md) across the Overview / zh-cn / de-de sheets, and appends a
# fourth data row (callerMd2.md) that didn't exist before.

$wb = $excel.ActiveWorkbook

$hlColor = 15570276  # RGB(0x64,0x95,0xED) == FF6495ED, matches the workbook's HyperLink style

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-22 19:04:41"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-22 19:04:41"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-22 19:04:41"

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-22 19:04:41"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null

$ws.Range("A2:A5").Font.Underline = $true
$ws.Range("A2:A5").Font.Color = $hlColor

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-22 19:04:37"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("J2").Value = "Include"
$ws.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-22 19:04:37"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"
$ws.Range("K3").Value = "e2e\callerMd1.md"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-22 19:04:37"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws.Range("J4").Value = "Include"

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-22 19:04:37"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "e2e\calleeMd1.md"
$ws.Range("J5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73ade9343d54c85e815177f52b3039f29ffd17e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73ade9343d54c85e815177f52b3039f29ffd17e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73ade9343d54c85e815177f52b3039f29ffd17e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73ade9343d54c85e815177f52b3039f29ffd17e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf") | Out-Null

$ws.Range("A2:A5").Font.Underline = $true
$ws.Range("A2:A5").Font.Color = $hlColor
$ws.Range("D2:D5").Font.Underline = $true
$ws.Range("D2:D5").Font.Color = $hlColor

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
$ws.Range("E2").Value = "2016-03-22 19:04:41"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("J2").Value = "Include"
$ws.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
$ws.Range("E3").Value = "2016-03-22 19:04:41"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"
$ws.Range("K3").Value = "e2e\callerMd1.md"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
$ws.Range("E4").Value = "2016-03-22 19:04:41"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws.Range("J4").Value = "Include"

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
$ws.Range("E5").Value = "2016-03-22 19:04:41"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "e2e\calleeMd1.md"
$ws.Range("J5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33ef7f8a00c2029321857dbdc078e9fa2fc41360/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33ef7f8a00c2029321857dbdc078e9fa2fc41360/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33ef7f8a00c2029321857dbdc078e9fa2fc41360/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/af94dd0d131653eb298d6096d937a2bb42da934b/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33ef7f8a00c2029321857dbdc078e9fa2fc41360/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf") | Out-Null

$ws.Range("A2:A5").Font.Underline = $true
$ws.Range("A2:A5").Font.Color = $hlColor
$ws.Range("D2:D5").Font.Underline = $true
$ws.Range("D2:D5").Font.Color = $hlColor

Write-Output "Report generated for handoff."
